# Applies the "Add files via upload" commit:
#  - Removes Sheet2 and Sheet3 (workbook now only contains Sheet1)
#  - Updates a handful of score values on Sheet1 (columns C/D) plus the
#    grand-total row (row 93) that rolls them up

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Sheet1")

# Remove the two blank extra sheets
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# Updated score values
$ws.Range("C2").Value = 26
$ws.Range("D7").Value = 57
$ws.Range("C14").Value = 97
$ws.Range("C28").Value = 57
$ws.Range("D28").Value = 47
$ws.Range("D34").Value = 43
$ws.Range("C36").Value = 43
$ws.Range("D36").Value = 27
$ws.Range("D37").Value = 43
$ws.Range("D46").Value = 56
$ws.Range("C53").Value = 84
$ws.Range("C57").Value = 71
$ws.Range("C61").Value = 40
$ws.Range("D61").Value = 34
$ws.Range("C70").Value = 54
$ws.Range("D70").Value = 44
$ws.Range("C78").Value = 81
$ws.Range("C81").Value = 76
$ws.Range("D81").Value = 54

# "任意登録" (optional-registration) subtotal row
$ws.Range("C92").Value = 253
$ws.Range("D92").Value = 187

# "総計" (grand total) row
$ws.Range("C93").Value = 5528
$ws.Range("D93").Value = 4459
